# Auto-generated edit script applying the committed numeric changes
# to the Ravana_Profits leve-profit sheets (scheduled-runner price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 917816.5600000001
$ws.Range("J100").Value = 10333.333
$ws.Range("L100").Value = 10333.333
$ws.Range("N100").Value = -11415.333
$ws.Range("H137").Value = 3665.6
$ws.Range("I137").Value = 1466.5
$ws.Range("J137").Value = 5131.6665
$ws.Range("K137").Value = 4399.5
$ws.Range("L137").Value = 15394.9995
$ws.Range("M137").Value = -1849.5
$ws.Range("N137").Value = -20494.9995
$ws.Range("H138").Value = 4716.273
$ws.Range("J138").Value = 4798
$ws.Range("L138").Value = 14394
$ws.Range("N138").Value = -24674
$ws.Range("H141").Value = 2000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 786.3333
$ws.Range("I5").Value = 783.8
$ws.Range("J5").Value = 799
$ws.Range("K5").Value = 783.8
$ws.Range("L5").Value = 799
$ws.Range("M5").Value = -671.8
$ws.Range("N5").Value = -1023
$ws.Range("H6").Value = 7333.6665
$ws.Range("I6").Value = 7333.6665
$ws.Range("K6").Value = 7333.6665
$ws.Range("M6").Value = -7160.6665
$ws.Range("H32").Value = 31293.75
$ws.Range("I32").Value = 31293.75
$ws.Range("K32").Value = 31293.75
$ws.Range("M32").Value = -31006.75
$ws.Range("H45").Value = 1885.8889
$ws.Range("J45").Value = 2483
$ws.Range("L45").Value = 2483
$ws.Range("N45").Value = -3237
$ws.Range("H74").Value = 3255
$ws.Range("I74").Value = 2257.1428
$ws.Range("J74").Value = 3792.3076
$ws.Range("K74").Value = 2257.1428
$ws.Range("L74").Value = 3792.3076
$ws.Range("M74").Value = -1383.1428
$ws.Range("N74").Value = -5540.3076
$ws.Range("H77").Value = 3255
$ws.Range("I77").Value = 2257.1428
$ws.Range("J77").Value = 3792.3076
$ws.Range("K77").Value = 11285.714
$ws.Range("L77").Value = 18961.538
$ws.Range("M77").Value = -6917.714
$ws.Range("N77").Value = -27697.538
$ws.Range("H97").Value = 237
$ws.Range("I97").Value = 237
$ws.Range("K97").Value = 237
$ws.Range("M97").Value = 259
$ws.Range("H102").Value = 2591.4546
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 2920.2666
$ws.Range("I132").Value = 2113.375
$ws.Range("K132").Value = 6340.125
$ws.Range("M132").Value = -3810.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 786.3333
$ws.Range("I4").Value = 783.8
$ws.Range("J4").Value = 799
$ws.Range("K4").Value = 783.8
$ws.Range("L4").Value = 799
$ws.Range("M4").Value = -668.8
$ws.Range("N4").Value = -1029
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H86").Value = 6133.615
$ws.Range("I86").Value = 5271
$ws.Range("J86").Value = 8074.5
$ws.Range("K86").Value = 5271
$ws.Range("L86").Value = 8074.5
$ws.Range("M86").Value = -4148
$ws.Range("N86").Value = -10320.5
$ws.Range("H89").Value = 6133.615
$ws.Range("I89").Value = 5271
$ws.Range("J89").Value = 8074.5
$ws.Range("K89").Value = 26355
$ws.Range("L89").Value = 40372.5
$ws.Range("M89").Value = -20739
$ws.Range("N89").Value = -51604.5
$ws.Range("H99").Value = 2001.4
$ws.Range("I99").Value = 2002
$ws.Range("J99").Value = 1999
$ws.Range("K99").Value = 2002
$ws.Range("L99").Value = 1999
$ws.Range("M99").Value = -504
$ws.Range("N99").Value = -4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 456.66666
$ws.Range("I7").Value = 493.45456
$ws.Range("J7").Value = 52
$ws.Range("K7").Value = 493.45456
$ws.Range("L7").Value = 52
$ws.Range("M7").Value = -380.45456
$ws.Range("N7").Value = -278
$ws.Range("H31").Value = 5645.3335
$ws.Range("I31").Value = 3258.5715
$ws.Range("J31").Value = 13999
$ws.Range("K31").Value = 3258.5715
$ws.Range("L31").Value = 13999
$ws.Range("M31").Value = -2963.5715
$ws.Range("N31").Value = -14589
$ws.Range("H34").Value = 5645.3335
$ws.Range("I34").Value = 3258.5715
$ws.Range("J34").Value = 13999
$ws.Range("K34").Value = 3258.5715
$ws.Range("L34").Value = 13999
$ws.Range("M34").Value = -3056.5715
$ws.Range("N34").Value = -14403
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H58").Value = 2471.0908
$ws.Range("I58").Value = 2471.0908
$ws.Range("K58").Value = 2471.0908
$ws.Range("M58").Value = -2268.0908
$ws.Range("H122").Value = 2003
$ws.Range("J122").Value = 3000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900
$ws.Range("H136").Value = 2471.0908
$ws.Range("I136").Value = 2471.0908
$ws.Range("K136").Value = 7413.2724
$ws.Range("M136").Value = -4863.2724
$ws.Range("H141").Value = 124898.5
$ws.Range("J141").Value = 133198
$ws.Range("L141").Value = 133198
$ws.Range("N141").Value = -143558

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 897
$ws.Range("I5").Value = 946.25
$ws.Range("K5").Value = 2838.75
$ws.Range("M5").Value = -2726.75
$ws.Range("H75").Value = 2416
$ws.Range("J75").Value = 2416
$ws.Range("L75").Value = 7248
$ws.Range("N75").Value = -9244
$ws.Range("H78").Value = 2416
$ws.Range("J78").Value = 2416
$ws.Range("L78").Value = 21744
$ws.Range("N78").Value = -31728
$ws.Range("H86").Value = 376
$ws.Range("I86").Value = 376
$ws.Range("K86").Value = 1128
$ws.Range("M86").Value = 58
$ws.Range("H87").Value = 2444.6667
$ws.Range("I87").Value = 2444.6667
$ws.Range("K87").Value = 7334.000100000001
$ws.Range("M87").Value = -6086.000100000001
$ws.Range("H89").Value = 376
$ws.Range("I89").Value = 376
$ws.Range("K89").Value = 3384
$ws.Range("M89").Value = 2544
$ws.Range("H90").Value = 2444.6667
$ws.Range("I90").Value = 2444.6667
$ws.Range("K90").Value = 22002.0003
$ws.Range("M90").Value = -15762.0003
$ws.Range("H117").Value = 1293.5
$ws.Range("J117").Value = 1250
$ws.Range("L117").Value = 3750
$ws.Range("N117").Value = -10634
$ws.Range("H127").Value = 1000
$ws.Range("J127").Value = 1000
$ws.Range("L127").Value = 3000
$ws.Range("N127").Value = -12920
$ws.Range("H135").Value = 897
$ws.Range("I135").Value = 946.25
$ws.Range("K135").Value = 8516.25
$ws.Range("M135").Value = -5981.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2220
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 2525
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 2525
$ws.Range("M43").Value = -849
$ws.Range("N43").Value = -2827
$ws.Range("H46").Value = 101490
$ws.Range("J46").Value = 101490
$ws.Range("L46").Value = 101490
$ws.Range("N46").Value = -101802
$ws.Range("H57").Value = 22000
$ws.Range("J57").Value = 22000
$ws.Range("L57").Value = 22000
$ws.Range("N57").Value = -23640
$ws.Range("H122").Value = 2148
$ws.Range("I122").Value = 1722
$ws.Range("K122").Value = 5166
$ws.Range("M122").Value = -2716
$ws.Range("H132").Value = 2842.7646
$ws.Range("I132").Value = 2486.923
$ws.Range("K132").Value = 7460.768999999999
$ws.Range("M132").Value = -4930.768999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 624994.25
$ws.Range("J43").Value = 624994.25
$ws.Range("L43").Value = 624994.25
$ws.Range("N43").Value = -625380.25
$ws.Range("H55").Value = 491
$ws.Range("I55").Value = 480
$ws.Range("J55").Value = 496.5
$ws.Range("K55").Value = 480
$ws.Range("L55").Value = 496.5
$ws.Range("M55").Value = -307
$ws.Range("N55").Value = -842.5
$ws.Range("H136").Value = 7177.1577
$ws.Range("I136").Value = 7147.9375
$ws.Range("K136").Value = 21443.8125
$ws.Range("M136").Value = -18893.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 3766.6667
$ws.Range("J54").Value = 3766.6667
$ws.Range("L54").Value = 3766.6667
$ws.Range("N54").Value = -4806.6667
$ws.Range("H132").Value = 3344.182
$ws.Range("I132").Value = 3818
$ws.Range("J132").Value = 2949.3333
$ws.Range("K132").Value = 11454
$ws.Range("L132").Value = 8847.999899999999
$ws.Range("M132").Value = -8924
$ws.Range("N132").Value = -13907.9999

Write-Output "Applied 237 cell updates across 8 sheets."
